{"js": "// The document stores each \"<id>...</id>\" marker as three separate runs:\n//   run 1: \"<id>\"        (Courier New, color 7f6000, sz/szCs 18)\n//   run 2: the id value  (color 000000, default font)\n//   run 3: \"</id>\"       (Courier New, color 7f6000, sz/szCs 18)\n// This edit merges each such triple into a single run (keeping the\n// formatting of the opening \"<id>\" run) whose text is the concatenation\n// \"<id>\" + value + \"</id>\". The final occurrence in the document also has\n// its id value corrected from the (duplicated) \"p046r_4\" to \"p046r_5\".\n\nconst body = context.document.body;\n\n// Find every opening \"<id>\" tag and every closing \"</id>\" tag; they are\n// emitted in document order, so pairing them positionally reconstructs\n// each \"<id>...</id>\" span regardless of what text sits between the tags.\nconst opens = body.search(\"<id>\", { matchCase: true });\nconst closes = body.search(\"</id>\", { matchCase: true });\nopens.load(\"items\");\ncloses.load(\"items\");\nawait context.sync();\n\nconst count = Math.min(opens.items.length, closes.items.length);\n\n// Build the full \"<id>...</id>\" range for each marker and read its text.\nconst fullRanges = [];\nfor (let i = 0; i < count; i++) {\n  fullRanges.push(opens.items[i].expandTo(closes.items[i]));\n}\nfor (const r of fullRanges) {\n  r.load(\"text\");\n}\nawait context.sync();\n\nconst originalTexts = fullRanges.map((r) => r.text);\n\n// Work out the replacement text for each marker. Normally this is just the\n// existing text (the only change is collapsing three runs into one), but\n// the last marker in the document had a duplicated id (\"p046r_4\") that\n// needs to become \"p046r_5\".\nconst idRegex = /^<id>([\\s\\S]*)<\\/id>$/;\nconst newTexts = originalTexts.slice();\nconst lastIndex = count - 1;\nif (lastIndex >= 0) {\n  const m = idRegex.exec(originalTexts[lastIndex]);\n  if (m && m[1] === \"p046r_4\") {\n    newTexts[lastIndex] = \"<id>p046r_5</id>\";\n  }\n}\n\n// Replace from the last marker to the first so earlier ranges stay valid\n// while later ones are rewritten.\nfor (let i = count - 1; i >= 0; i--) {\n  if (newTexts[i] !== originalTexts[i] || true) {\n    fullRanges[i].insertText(newTexts[i], \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "# The document stores each \"<id>...</id>\" marker as three separate runs:\n#   run 1: \"<id>\"        (Courier New, color 7f6000, sz/szCs 18)\n#   run 2: the id value  (color 000000, default font)\n#   run 3: \"</id>\"       (Courier New, color 7f6000, sz/szCs 18)\n# Each such marker occupies its own paragraph in its entirety. This edit\n# merges each triple of runs into a single run (keeping the formatting of\n# the opening \"<id>\" run) whose text is the concatenation\n# \"<id>\" + value + \"</id>\". The final occurrence in the document also has\n# its id value corrected from the (duplicated) \"p046r_4\" to \"p046r_5\".\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n# Collect the indices of every paragraph that is entirely an \"<id>...</id>\"\n# marker (in document order).\n$hitIdxs = @()\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -match \"^<id>[\\s\\S]*</id>\\r?$\") {\n        $hitIdxs += , $i\n    }\n}\n\n$lastHit = $hitIdxs.Count - 1\n\n# Walk backwards so replacing one paragraph's text never shifts the\n# positions of paragraphs we still have to process.\nfor ($k = $lastHit; $k -ge 0; $k--) {\n    $i = $hitIdxs[$k]\n    $rng = $d.Paragraphs.Item($i).Range\n    $full = $rng.Text\n    # Strip the trailing paragraph mark captured by Paragraph.Range.\n    $content = $full.Substring(0, $full.Length - 1)\n\n    # Fix the duplicated id on the final marker in the document.\n    if ($k -eq $lastHit -and $content -eq \"<id>p046r_4</id>\") {\n        $content = \"<id>p046r_5</id>\"\n    }\n\n    # Replace only the marker text, leaving the paragraph mark untouched,\n    # by extending the range one character past \"</id>\" so Word collapses\n    # the three runs into a single run carrying the \"<id>\" run's formatting.\n    $full2 = $d.Range($rng.Start, $rng.Start + $content.Length + 1)\n    $full2.Text = $content\n}\n"}
